# Add a new column C: a header "date stamp" cell (styled like the other
# header cells) plus per-row Present/Absent markers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- C1 header cell -------------------------------------------------
# We want C1 to end up formatted exactly like the existing header cells
# (A1/B1: bold font, thin border all around, centered/top-aligned) while
# holding the literal text "2025-07-03" (NOT an auto-converted date
# serial). Force text storage first (NumberFormat "@") so Excel doesn't
# coerce the date-looking string into a date value, then copy the
# header formatting from B1 on top so the final cell style matches the
# other header cells.
$ws.Range("C1").NumberFormat = "@"
$ws.Range("C1").Value = "2025-07-03"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats

# --- C2:C6 attendance markers ---------------------------------------
$ws.Range("C2").Value = "P"
$ws.Range("C3").Value = "A"
$ws.Range("C4").Value = "A"
$ws.Range("C5").Value = "A"
$ws.Range("C6").Value = "A"
